# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp (10:35 -> 11:05)
# - Refresh case counts for Estados Unidos, Indonesia, Austria, Oman
# - Refresh case counts for Banglades, which now overtakes Emiratos Arabes
#   Unidos in the ranking, so the two rows swap places (row 32 <-> row 33)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "last refreshed" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 11:05"

# Estados Unidos (row 4) - updated totals
$ws.Range("B4").Value = 1527935
$ws.Range("C4").Value = 271
$ws.Range("E4").Value = 1090568

# Banglades moves up to row 32 (was row 33) with refreshed numbers,
# overtaking Emiratos Arabes Unidos which drops to row 33 unchanged.
$ws.Range("A32").Value = "Banglades"
$ws.Range("B32").Value = 23870
$ws.Range("C32").Value = 1602
$ws.Range("D32").Value = 4585
$ws.Range("E32").Value = 18936
$ws.Range("G32").Value = 21
$ws.Range("H32").Value = 349

$ws.Range("A33").Value = "Emiratos Arabes Unidos"
$ws.Range("B33").Value = 23358
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 8512
$ws.Range("E33").Value = 14626
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 220

# Indonesia (row 36) - updated totals
$ws.Range("B36").Value = 18010
$ws.Range("C36").Value = 496
$ws.Range("D36").Value = 4324
$ws.Range("E36").Value = 12495
$ws.Range("G36").Value = 43
$ws.Range("H36").Value = 1191

# Austria (row 40) - updated totals
$ws.Range("B40").Value = 16269
$ws.Range("C40").Value = 27
$ws.Range("D40").Value = 14614
$ws.Range("E40").Value = 1026

# Oman (row 65) - updated totals
$ws.Range("B65").Value = 5379
$ws.Range("C65").Value = 193
$ws.Range("D65").Value = 1496
$ws.Range("E65").Value = 3860
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = 23
